# Refresh the chi-square / feature-importance statistics on the "factors" sheet
# (recomputed after switching the incident data loader to read CSV files),
# and swap the row order of the self_service / has_knowledge_article factors
# to match the new column numbering.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"5086.153335856128"
$ws.Range("H2").Value = [double]"0.8540265712417093"
$ws.Range("I2").Value = [double]"0.04997149770573032"

$ws.Range("F3").Value = [double]"1204.914932712188"
$ws.Range("G3").Value = [double]"1.439258204408449e-247"
$ws.Range("H3").Value = [double]"0.0675069979266494"
$ws.Range("I3").Value = [double]"0.02146796566642512"

$ws.Range("F4").Value = [double]"220.8241972240276"
$ws.Range("G4").Value = [double]"5.978698520033615e-50"
$ws.Range("H4").Value = [double]"0.03113305331997251"
$ws.Range("I4").Value = [double]"0.01358754557154232"

$ws.Range("F5").Value = [double]"0.8883290064550426"
$ws.Range("G5").Value = [double]"0.3459305333140232"
$ws.Range("H5").Value = [double]"0.00492061072439518"
$ws.Range("I5").Value = [double]"0.002660039725387586"

$ws.Range("F6").Value = [double]"86.00805738538253"
$ws.Range("G6").Value = [double]"1.792039823210119e-20"
$ws.Range("H6").Value = [double]"0.002000261713705436"
$ws.Range("I6").Value = [double]"0.002567947345023552"

$ws.Range("F7").Value = [double]"392.3837685207457"
$ws.Range("G7").Value = [double]"2.752934587762642e-74"
$ws.Range("H7").Value = [double]"0.01301928421089169"
$ws.Range("I7").Value = [double]"0.00237484453601966"

$ws.Range("A8").Value = "self_service"
$ws.Range("B8").Value = [double]"8"
$ws.Range("F8").Value = [double]"42.314501917315"
$ws.Range("G8").Value = [double]"7.771416267429973e-11"
$ws.Range("H8").Value = [double]"0.0003918017624877692"
$ws.Range("I8").Value = [double]"0.0007551729345986052"

$ws.Range("F9").Value = [double]"116.6857017813586"
$ws.Range("G9").Value = [double]"3.363625883314413e-27"
$ws.Range("H9").Value = [double]"0.02110987240302388"

$ws.Range("F10").Value = [double]"915.4726472302203"
$ws.Range("G10").Value = [double]"4.248980687842385e-201"
$ws.Range("H10").Value = [double]"0.003545907502174473"

$ws.Range("A11").Value = "has_knowledge_article"
$ws.Range("B11").Value = [double]"9"
$ws.Range("F11").Value = [double]"34.81280277953115"
$ws.Range("G11").Value = [double]"3.629786486834646e-09"
$ws.Range("H11").Value = [double]"0.002345639194990277"

$ws.Range("F12").Value = [double]"109.391853707005"
$ws.Range("G12").Value = [double]"1.331757529725758e-25"
$ws.Range("H12").Value = [double]"0"

$ws.Range("F13").Value = [double]"30.0312876065441"
$ws.Range("G13").Value = [double]"4.251311960370145e-08"

$ws.Range("F14").Value = [double]"21.56381917347364"
$ws.Range("G14").Value = [double]"3.422476790024476e-06"

$ws.Range("F15").Value = [double]"15.28524034640144"
$ws.Range("G15").Value = [double]"9.243595384270451e-05"

$ws.Range("F16").Value = [double]"3.540962606063545"
$ws.Range("G16").Value = [double]"0.05987070582593627"

$ws.Range("E17").Value = [double]"1988"
$ws.Range("F17").Value = [double]"2681.182401149468"
$ws.Range("G17").Value = [double]"1.250721961446571e-23"

$ws.Range("E18").Value = [double]"838"
$ws.Range("F18").Value = [double]"1824.123930449769"
$ws.Range("G18").Value = [double]"2.858346425893734e-75"

$ws.Range("E19").Value = [double]"118"
$ws.Range("F19").Value = [double]"757.5341707421819"
$ws.Range("G19").Value = [double]"2.201367192346436e-94"

$ws.Range("F20").Value = [double]"335.8134709170386"
$ws.Range("G20").Value = [double]"6.459790618508831e-67"

$ws.Range("F21").Value = [double]"103.315517062819"
$ws.Range("G21").Value = [double]"1.180107952614957e-17"

$ws.Range("F22").Value = [double]"75.24537910926617"
$ws.Range("G22").Value = [double]"1.76813754184472e-15"

$ws.Range("F23").Value = [double]"53.02979898947834"
$ws.Range("G23").Value = [double]"7.351678343089048e-08"

$ws.Range("F24").Value = [double]"23.80406556053098"
$ws.Range("G24").Value = [double]"2.744679293689399e-05"

$ws.Range("F25").Value = [double]"0.2463471129800019"
$ws.Range("G25").Value = [double]"0.6196590139590284"
